$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.088.41'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.43%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.600.52'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.09%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.13'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.91%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3779'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.62%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3649'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.38%  '

# Row 9
$ws.Range("E9").Value = '  -4.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.268'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.13%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08160'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.54%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.05'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.69%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.588'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.90%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001258'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.364'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -8.45%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.600.71'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.60'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06860'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.74%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.562'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.62%  '

# Row 22
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5555'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.78%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.96'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.81%  '

# Row 25
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '23.086.06'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.40%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.342'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.47%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.715'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -8.50%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.13'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.34%  '

# Row 29
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.37'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.03%  '

# Row 30
$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.273'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.69%  '

# Row 31
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.39'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.08%  '

# Row 32
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.402'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.15%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.839'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -12.99%  '

# Row 34
$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.776.19'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.26%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9621'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.53%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07651'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.52%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.255'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.99%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02725'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.89%  '

# Row 39
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2549'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.89%  '

# Row 40
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08905'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.41%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.06'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.17%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.369'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.87%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7091'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -6.49%  '

# Row 44
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.66'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.11%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.46'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.81%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6601'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.94%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.05%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.309'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.27%  '

# Row 49
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.984'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.68%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.97'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.88%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07935'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.38%  '
